# "Pasar divisor a constantes y leer propiedades de Header de archivo"
#
# The worksheet's header row (row 1) renamed its two "Out Col N" labels to
# the shorter "OutputN" form. Columns A and C are best-fit to their
# contents, so shortening the header text also shrinks their width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two header cells (A1 = "Out Col 7" -> "Output7",
# C1 = "Out Col 1" -> "Output1").
$ws.Range("A1").Value = "Output7"
$ws.Range("C1").Value = "Output1"

# Re-apply best-fit width on the two header columns now that their
# (shorter) text changed, shrinking them from ~10.43 to ~9.43 characters.
$ws.Columns.Item(1).ColumnWidth = 8.666666666666668
$ws.Columns.Item(3).ColumnWidth = 8.666666666666668
